$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: replace the old "M2.5 x 8 Umbraco bolt" component with the new
# "M3 x 8 Umbraco bolt" model used for the encoder knob.
$ws.Range("A6").Value = "M3 x 8 Umbraco bolt"

# E10 currently shows the AliExpress URL as plain text with no live
# hyperlink target (a "latch" solenoid lock). Give it a hyperlink
# display annotation (Excel caps TextToDisplay at 255 characters).
$full = $ws.Range("E10").Text
$display = $full.Substring(0, 255)
$ws.Hyperlinks.Add($ws.Range("E10"), "", "", "", $display) | Out-Null

# Move the active selection to C7, matching the saved view state.
$ws.Range("C7").Select()
